# Apply the edit described in the diff to kayser_flow_time.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-driven values in B2:B7 with new literal (computed) values
$ws.Range("B2").Value = 0.0131862063225597
$ws.Range("B3").Value = 0.000025730435380530201
$ws.Range("B4").Value = 0.000000090282457069618899
$ws.Range("B5").Value = 0.000000084933130121127702
$ws.Range("B6").Value = 0.000000073104273560166196
$ws.Range("B7").Value = 0.00000000042192706892422999

# Clear the remaining cycNo (A) entries and cssCrit (B) formulas for rows 8-16,
# leaving the cells blank (but still formatted)
$ws.Range("A8:A16").ClearContents()
$ws.Range("B8:B16").ClearContents()

# Remove rows 22-30 entirely, shrinking the sheet's used range down to row 21
$ws.Range("A22:B30").EntireRow.Delete()

# Update the worksheet view's selected cell/range
$ws.Range("E5").Select()
